$wb = $excel.ActiveWorkbook

# --- Sheet "RUN_MANAGER" (sheet1) ---
$ws1 = $wb.Worksheets.Item("RUN_MANAGER")

# C3: execute yes -> no  (newTest row)
$ws1.Range("C3").Value = "no"

# New row 4: amazonPageTest
$ws1.Range("A4").Value = "amazonPageTest"
$ws1.Range("B4").Value = "Amazon-Page-Test"
$ws1.Range("C4").Value = "yes"
$ws1.Range("D4").Value = "'1"
$ws1.Range("E4").Value = "'1"

# D2: count 1 -> 2  (loginLogoutTest row) -- added last so the shared
# string "2" is appended after the new unique strings above
$ws1.Range("D2").Value = "'2"

$ws1.Range("D12").Select()

# --- Sheet "DATA" (sheet2) ---
$ws2 = $wb.Worksheets.Item("DATA")

# B3: execute yes -> no  (loginLogoutTest/firefox row)
$ws2.Range("B3").Value = "no"

# B4: execute yes -> no  (newTest/firefox row)
$ws2.Range("B4").Value = "no"

# B6: execute no -> yes  (loginLogoutTest/chrome row)
$ws2.Range("B6").Value = "yes"

# New row 7: amazonPageTest
$ws2.Range("A7").Value = "amazonPageTest"
$ws2.Range("B7").Value = "yes"
$ws2.Range("C7").Value = "chrome"
$ws2.Range("D7").Value = "'"
$ws2.Range("E7").Value = "'"
$ws2.Range("F7").Value = "'"

$ws2.Range("B7").Select()
